# Update the Outflow Q and WQ worksheet (pMeHg_Comb) for wy2016-17:
#  - collapse the separate "PpMeHg"/"RpMeHg" (particulate/remark) columns into
#    a single "pMeHg" column, dropping the now-unused qualifier column (E)
#  - revise two model results (rows 4 & 5) from 0.1 to 0.05 now that the
#    "<" remark-code column is gone, keeping the highlighted fill but right
#    aligning the values
#  - give row 15's pMeHg value the same 2-decimal number format as the rest
#    of the column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The RpMeHg column (E) is removed entirely; its header/2-letter-code/data
# go away and everything to its right (nothing, here) shifts left.
$ws.Columns("E").Delete()

# Header: "PpMeHg" -> "pMeHg"
$ws.Range("D1").Value = "pMeHg"

# Revised pMeHg results for 1/23/2016 and 1/24/2016
$ws.Range("D4").Value = 0.05
$ws.Range("D5").Value = 0.05

# Keep the highlighted fill on D4:D5 but now right-align like the other
# numeric/date/time columns
$ws.Range("D4:D5").HorizontalAlignment = -4152   # xlRight

# D15 now carries an explicit 2-decimal number format
$ws.Range("D15").NumberFormat = "0.00"

# Matches the author's final selection after the edit (one cell right of the
# new last column, same relative spot as before the column delete)
[void]$ws.Range("E31").Select()
